# Applies: added JW and Bible Text
# Adds four new verse rows (7-10) to the "Mathew" sheet, with Adhola text in
# column A and chapter/verse numbers in column C, then updates the sheet
# selection to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mathew")

# Row 7
$ws.Range("A7").Value2 = "Ci Kerode olwoŋo luryeko i muŋ, guniaŋe kun gititte kare ma lakalatwe onen iye"
$ws.Range("C7").Value2 = 2.7
$ws.Rows.Item(7).RowHeight = 30

# Row 8
$ws.Range("A8").Value2 = "Ka doŋ ocwalogi me cito i Jerucalem kun waco botgi ni, “Wuciti wupeny lok kom latin man maber kikore, ce ka wunoŋe ci wudwokka lok, wek an bene acit awore.”"
$ws.Range("C8").Value2 = 2.8
$ws.Rows.Item(8).RowHeight = 45

# Row 9
$ws.Range("A9").Value2 = "I kare ma doŋ guwinyo lok pa kabaka, gucito ki yogi ci guneno lakalatwe ma yam koŋ guneno yo tuŋ wokceŋ ca otelo nyimgi, obino ocuŋ ki malo wa i kabedo ka ma onoŋo latin-nu tye iye"
$ws.Range("C9").Value2 = 2.9
$ws.Rows.Item(9).RowHeight = 45

# Row 10
$ws.Range("A10").Value2 = "Ka guneno lakalatwe meno, cwinygi obedo yom adada."

# Update the selected cell shown in the workbook to A4
$ws.Range("A4").Select()
